$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.392.21'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.434.93'
$ws.Range("E3").Value = '  +1.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.77'
$ws.Range("E5").Value = '  -1.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.47'
$ws.Range("E6").Value = '  -1.93%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.427.90'
$ws.Range("E7").Value = '  +1.23%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("E9").Value = '  -0.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.198'
$ws.Range("E10").Value = '  +0.70%  '

$ws.Range("E11").Value = '  -0.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.81'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("E13").Value = '  -1.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '699.26'
$ws.Range("E14").Value = '  +1.79%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.65'
$ws.Range("E15").Value = '  +0.53%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.962.76'
$ws.Range("E16").Value = '  +0.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.369.22'
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.431.73'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("E19").Value = '  +0.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.74'
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.900'
$ws.Range("E22").Value = '  -0.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.40'
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("E24").Value = '  -0.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '101.16'
$ws.Range("E25").Value = '  -3.21%  '

$ws.Range("E26").Value = '  -1.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.68'
$ws.Range("E27").Value = '  -1.77%  '

$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.63'
$ws.Range("E29").Value = '  -2.82%  '

$ws.Range("E30").Value = '  +0.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.98'
$ws.Range("E31").Value = '  -0.91%  '

$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.76'
$ws.Range("E32").Value = '  +2.54%  '

$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '566.91'
$ws.Range("E33").Value = '  +1.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.04'
$ws.Range("E34").Value = '  -1.29%  '

$ws.Range("E35").Value = '  -1.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.16'
$ws.Range("E36").Value = '  -0.38%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.610.17'
$ws.Range("E38").Value = '  -3.19%  '

$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.99'
$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0732'
$ws.Range("E41").Value = '  +3.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.29'
$ws.Range("E42").Value = '  +1.32%  '

$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.35'
$ws.Range("E44").Value = '  +3.18%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.333'
$ws.Range("E45").Value = '  -2.02%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0420'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.47'
$ws.Range("E47").Value = '  +4.57%  '

$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.66'
$ws.Range("E48").Value = '  +0.05%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.128'
$ws.Range("E49").Value = '  -1.26%  '

$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.30'
$ws.Range("E51").Value = '  -1.09%  '
